$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest crypto price/volume snapshot refresh (GitHub Actions scheduled run).
# Maps worksheet row -> updated Price (column D) / Volume(1h) (column E) text.
# A value of $null means that column did not change for that row.
$updates = @{
    2 = @{ D = '30.318.47'; E = '  +0.24%  ' };
    3 = @{ D = '1.869.83'; E = '  +0.30%  ' };
    4 = @{ D = $null; E = $null };
    5 = @{ D = '235.22'; E = '  -0.81%  ' };
    6 = @{ D = $null; E = '  +0.06%  ' };
    7 = @{ D = '0.4696'; E = '  +0.31%  ' };
    8 = @{ D = $null; E = '  -0.19%  ' };
    9 = @{ D = '0.06598'; E = '  +0.69%  ' };
    10 = @{ D = '21.80'; E = '  -1.75%  ' };
    11 = @{ D = '0.07985'; E = '  +1.18%  ' };
    12 = @{ D = '96.85'; E = '  -1.15%  ' };
    13 = @{ D = '1.874.21'; E = '  +0.41%  ' };
    14 = @{ D = '0.6897'; E = '  +0.92%  ' };
    15 = @{ D = $null; E = '  -1.56%  ' };
    16 = @{ D = '269.07'; E = '  -3.16%  ' };
    17 = @{ D = '30.353.43'; E = '  +0.35%  ' };
    18 = @{ D = '14.14'; E = '  +3.88%  ' };
    19 = @{ D = '0.000007791'; E = '  +6.01%  ' };
    20 = @{ D = $null; E = '  +0.06%  ' };
    21 = @{ D = '2.119.74'; E = '  +0.03%  ' };
    22 = @{ D = '1.001'; E = '  +0.00%  ' };
    23 = @{ D = '5.262'; E = '  -1.85%  ' };
    24 = @{ D = '6.217'; E = '  +0.29%  ' };
    25 = @{ D = '9.396'; E = '  +1.58%  ' };
    26 = @{ D = '167.53'; E = '  -0.47%  ' };
    27 = @{ D = $null; E = '  -1.01%  ' };
    28 = @{ D = '1.950'; E = $null };
    29 = @{ D = '1.365'; E = $null };
    30 = @{ D = '0.09864'; E = '  +0.16%  ' };
    31 = @{ D = $null; E = '  -0.76%  ' };
    32 = @{ D = '1.459'; E = '  -1.59%  ' };
    33 = @{ D = $null; E = '  -0.44%  ' };
    34 = @{ D = '0.04712'; E = '  -0.76%  ' };
    35 = @{ D = '1.137'; E = '  -0.09%  ' };
    36 = @{ D = '0.7026'; E = '  -0.35%  ' };
    37 = @{ D = '2.738'; E = '  +1.13%  ' };
    38 = @{ D = $null; E = '  -0.05%  ' };
    39 = @{ D = '2.811'; E = '  +6.81%  ' };
    40 = @{ D = $null; E = '  -0.63%  ' };
    41 = @{ D = '72.17'; E = '  -4.49%  ' };
    42 = @{ D = $null; E = '  -0.02%  ' };
    43 = @{ D = '0.4178'; E = '  -0.10%  ' };
    44 = @{ D = '0.8426'; E = '  -1.30%  ' };
    45 = @{ D = '1.001'; E = '  +0.07%  ' };
    46 = @{ D = '103.35'; E = '  -0.21%  ' };
    47 = @{ D = '7.093'; E = '  -1.72%  ' };
    48 = @{ D = '9.113'; E = '  -1.64%  ' };
    49 = @{ D = '919.27'; E = '  -3.16%  ' };
    50 = @{ D = '34.52'; E = '  +0.77%  ' };
    51 = @{ D = $null; E = '  +0.91%  ' }
}

# Rows whose new Price text is "pure-numeric-looking" (e.g. trailing zero
# like '21.80' or '1.950') need the cell pre-formatted as Text; otherwise
# Excel's normal numeric auto-detection on assignment would silently drop
# the significant trailing zero (21.80 -> 21.8), which would not match the
# source price feed's string formatting.
$forceTextRows = @(10, 28)

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.D) {
        $cell = $ws.Cells.Item($row, 4)
        if ($forceTextRows -contains $row) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $vals.D
    }
    if ($null -ne $vals.E) {
        $ws.Cells.Item($row, 5).Value = $vals.E
    }
}
